$d = $word.ActiveDocument

$d.Content.Find.Execute("77-22=55", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=13", 2) | Out-Null
$d.Content.Find.Execute("2+33=35", $true, $false, $false, $false, $false, $true, 1, $false, "95-72=23", 2) | Out-Null
$d.Content.Find.Execute("83-32=51", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=30", 2) | Out-Null
$d.Content.Find.Execute("80-9=71", $true, $false, $false, $false, $false, $true, 1, $false, "11-4=7", 2) | Out-Null
$d.Content.Find.Execute("90-15=75", $true, $false, $false, $false, $false, $true, 1, $false, "6+44=50", 2) | Out-Null
$d.Content.Find.Execute("5+29=34", $true, $false, $false, $false, $false, $true, 1, $false, "98-0=98", 2) | Out-Null
$d.Content.Find.Execute("33-30=3", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=37", 2) | Out-Null
$d.Content.Find.Execute("18+53=71", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("97-8=89", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=55", 2) | Out-Null
$d.Content.Find.Execute("0+42=42", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=82", 2) | Out-Null
$d.Content.Find.Execute("92-21=71", $true, $false, $false, $false, $false, $true, 1, $false, "6+20=26", 2) | Out-Null
$d.Content.Find.Execute("88-81=7", $true, $false, $false, $false, $false, $true, 1, $false, "69-27=42", 2) | Out-Null
$d.Content.Find.Execute("72-16=56", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=66", 2) | Out-Null
$d.Content.Find.Execute("1+34=35", $true, $false, $false, $false, $false, $true, 1, $false, "65-32=33", 2) | Out-Null
$d.Content.Find.Execute("83-36=47", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=22", 2) | Out-Null
$d.Content.Find.Execute("3+39=42", $true, $false, $false, $false, $false, $true, 1, $false, "10-0=10", 2) | Out-Null
$d.Content.Find.Execute("22+2=24", $true, $false, $false, $false, $false, $true, 1, $false, "85+9=94", 2) | Out-Null
$d.Content.Find.Execute("71-32=39", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("20+48=68", $true, $false, $false, $false, $false, $true, 1, $false, "21-6=15", 2) | Out-Null
$d.Content.Find.Execute("78-65=13", $true, $false, $false, $false, $false, $true, 1, $false, "15-0=15", 2) | Out-Null
$d.Content.Find.Execute("6+38=44", $true, $false, $false, $false, $false, $true, 1, $false, "86-84=2", 2) | Out-Null
$d.Content.Find.Execute("66-49=17", $true, $false, $false, $false, $false, $true, 1, $false, "85-73=12", 2) | Out-Null
$d.Content.Find.Execute("72+23=95", $true, $false, $false, $false, $false, $true, 1, $false, "0+8=8", 2) | Out-Null
$d.Content.Find.Execute("32+41=73", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=21", 2) | Out-Null
$d.Content.Find.Execute("81-35=46", $true, $false, $false, $false, $false, $true, 1, $false, "30-16=14", 2) | Out-Null
$d.Content.Find.Execute("11+72=83", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=51", 2) | Out-Null
$d.Content.Find.Execute("23-12=11", $true, $false, $false, $false, $false, $true, 1, $false, "43-17=26", 2) | Out-Null
$d.Content.Find.Execute("25+51=76", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=37", 2) | Out-Null
$d.Content.Find.Execute("87-5=82", $true, $false, $false, $false, $false, $true, 1, $false, "58-16=42", 2) | Out-Null
$d.Content.Find.Execute("76-26=50", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=16", 2) | Out-Null
$d.Content.Find.Execute("47-31=16", $true, $false, $false, $false, $false, $true, 1, $false, "75-60=15", 2) | Out-Null
$d.Content.Find.Execute("74-57=17", $true, $false, $false, $false, $false, $true, 1, $false, "51+43=94", 2) | Out-Null
$d.Content.Find.Execute("27-10=17", $true, $false, $false, $false, $false, $true, 1, $false, "65-32=33", 2) | Out-Null
$d.Content.Find.Execute("62-40=22", $true, $false, $false, $false, $false, $true, 1, $false, "68-67=1", 2) | Out-Null
$d.Content.Find.Execute("45+33=78", $true, $false, $false, $false, $false, $true, 1, $false, "79-45=34", 2) | Out-Null
$d.Content.Find.Execute("78-7=71", $true, $false, $false, $false, $false, $true, 1, $false, "68-48=20", 2) | Out-Null
$d.Content.Find.Execute("25+39=64", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=55", 2) | Out-Null
$d.Content.Find.Execute("2-1=1", $true, $false, $false, $false, $false, $true, 1, $false, "44-34=10", 2) | Out-Null
$d.Content.Find.Execute("9+82=91", $true, $false, $false, $false, $false, $true, 1, $false, "96-66=30", 2) | Out-Null
$d.Content.Find.Execute("89-76=13", $true, $false, $false, $false, $false, $true, 1, $false, "10+39=49", 2) | Out-Null
$d.Content.Find.Execute("59+38=97", $true, $false, $false, $false, $false, $true, 1, $false, "96-77=19", 2) | Out-Null
$d.Content.Find.Execute("7+1=8", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=43", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "45-14=31", 2) | Out-Null
$d.Content.Find.Execute("8+42=50", $true, $false, $false, $false, $false, $true, 1, $false, "77-2=75", 2) | Out-Null
$d.Content.Find.Execute("52-26=26", $true, $false, $false, $false, $false, $true, 1, $false, "82+1=83", 2) | Out-Null
$d.Content.Find.Execute("24+3=27", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=58", 2) | Out-Null
$d.Content.Find.Execute("52-5=47", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=75", 2) | Out-Null
$d.Content.Find.Execute("52-7=45", $true, $false, $false, $false, $false, $true, 1, $false, "60-52=8", 2) | Out-Null
$d.Content.Find.Execute("93-0=93", $true, $false, $false, $false, $false, $true, 1, $false, "25+13=38", 2) | Out-Null
$d.Content.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "65+32=97", 2) | Out-Null
$d.Content.Find.Execute("50-24=26", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=9", 2) | Out-Null
$d.Content.Find.Execute("58-0=58", $true, $false, $false, $false, $false, $true, 1, $false, "27+20=47", 2) | Out-Null
$d.Content.Find.Execute("48+34=82", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=26", 2) | Out-Null
$d.Content.Find.Execute("90-29=61", $true, $false, $false, $false, $false, $true, 1, $false, "4+3=7", 2) | Out-Null
$d.Content.Find.Execute("2+7=9", $true, $false, $false, $false, $false, $true, 1, $false, "54+2=56", 2) | Out-Null
$d.Content.Find.Execute("57+14=71", $true, $false, $false, $false, $false, $true, 1, $false, "61+35=96", 2) | Out-Null
$d.Content.Find.Execute("74-23=51", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=41", 2) | Out-Null
$d.Content.Find.Execute("63+3=66", $true, $false, $false, $false, $false, $true, 1, $false, "25+31=56", 2) | Out-Null
$d.Content.Find.Execute("39+0=39", $true, $false, $false, $false, $false, $true, 1, $false, "10+88=98", 2) | Out-Null
$d.Content.Find.Execute("22+73=95", $true, $false, $false, $false, $false, $true, 1, $false, "70+9=79", 2) | Out-Null
$d.Content.Find.Execute("83-78=5", $true, $false, $false, $false, $false, $true, 1, $false, "39+54=93", 2) | Out-Null
$d.Content.Find.Execute("1+91=92", $true, $false, $false, $false, $false, $true, 1, $false, "66-22=44", 2) | Out-Null
$d.Content.Find.Execute("32+11=43", $true, $false, $false, $false, $false, $true, 1, $false, "10+45=55", 2) | Out-Null
$d.Content.Find.Execute("8+11=19", $true, $false, $false, $false, $false, $true, 1, $false, "87-81=6", 2) | Out-Null
$d.Content.Find.Execute("3+13=16", $true, $false, $false, $false, $false, $true, 1, $false, "64-17=47", 2) | Out-Null
$d.Content.Find.Execute("78+16=94", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=35", 2) | Out-Null
$d.Content.Find.Execute("26-4=22", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=29", 2) | Out-Null
$d.Content.Find.Execute("70+24=94", $true, $false, $false, $false, $false, $true, 1, $false, "51+38=89", 2) | Out-Null
$d.Content.Find.Execute("99-5=94", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=11", 2) | Out-Null
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "6+77=83", 2) | Out-Null
$d.Content.Find.Execute("49+0=49", $true, $false, $false, $false, $false, $true, 1, $false, "53+21=74", 2) | Out-Null
$d.Content.Find.Execute("64-30=34", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("21-13=8", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=68", 2) | Out-Null
$d.Content.Find.Execute("45+12=57", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("19+25=44", $true, $false, $false, $false, $false, $true, 1, $false, "82-35=47", 2) | Out-Null
$d.Content.Find.Execute("72+22=94", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=59", 2) | Out-Null
$d.Content.Find.Execute("1+65=66", $true, $false, $false, $false, $false, $true, 1, $false, "73-4=69", 2) | Out-Null
$d.Content.Find.Execute("96-58=38", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("35+57=92", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2) | Out-Null
$d.Content.Find.Execute("79-1=78", $true, $false, $false, $false, $false, $true, 1, $false, "25-15=10", 2) | Out-Null
$d.Content.Find.Execute("6+76=82", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=1", 2) | Out-Null
$d.Content.Find.Execute("91-13=78", $true, $false, $false, $false, $false, $true, 1, $false, "54+0=54", 2) | Out-Null
$d.Content.Find.Execute("23-3=20", $true, $false, $false, $false, $false, $true, 1, $false, "28+31=59", 2) | Out-Null
$d.Content.Find.Execute("82-30=52", $true, $false, $false, $false, $false, $true, 1, $false, "60-10=50", 2) | Out-Null
$d.Content.Find.Execute("79-52=27", $true, $false, $false, $false, $false, $true, 1, $false, "6+7=13", 2) | Out-Null
$d.Content.Find.Execute("53-27=26", $true, $false, $false, $false, $false, $true, 1, $false, "89-11=78", 2) | Out-Null
$d.Content.Find.Execute("43-21=22", $true, $false, $false, $false, $false, $true, 1, $false, "72-55=17", 2) | Out-Null
$d.Content.Find.Execute("64-19=45", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=31", 2) | Out-Null
$d.Content.Find.Execute("16+41=57", $true, $false, $false, $false, $false, $true, 1, $false, "64-21=43", 2) | Out-Null
$d.Content.Find.Execute("77-71=6", $true, $false, $false, $false, $false, $true, 1, $false, "94-75=19", 2) | Out-Null
$d.Content.Find.Execute("59-37=22", $true, $false, $false, $false, $false, $true, 1, $false, "27+33=60", 2) | Out-Null
$d.Content.Find.Execute("35+45=80", $true, $false, $false, $false, $false, $true, 1, $false, "43-37=6", 2) | Out-Null
$d.Content.Find.Execute("66-15=51", $true, $false, $false, $false, $false, $true, 1, $false, "69-17=52", 2) | Out-Null
$d.Content.Find.Execute("2+30=32", $true, $false, $false, $false, $false, $true, 1, $false, "34-7=27", 2) | Out-Null
$d.Content.Find.Execute("55+27=82", $true, $false, $false, $false, $false, $true, 1, $false, "2+81=83", 2) | Out-Null
$d.Content.Find.Execute("31+38=69", $true, $false, $false, $false, $false, $true, 1, $false, "65-3=62", 2) | Out-Null
$d.Content.Find.Execute("77-20=57", $true, $false, $false, $false, $false, $true, 1, $false, "37+27=64", 2) | Out-Null
$d.Content.Find.Execute("52+3=55", $true, $false, $false, $false, $false, $true, 1, $false, "41+32=73", 2) | Out-Null
$d.Content.Find.Execute("23+41=64", $true, $false, $false, $false, $false, $true, 1, $false, "7-2=5", 2) | Out-Null
$d.Content.Find.Execute("43+26=69", $true, $false, $false, $false, $false, $true, 1, $false, "74+0=74", 2) | Out-Null
